# Apply the "Confirmations" worksheet relabeling + summary-row cleanup.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Value" -> "Values" -------------------------------------
$ws.Range("B1").Value = "Values"

# --- Civilian section (rows 7-11): prefix each sub-label with the
#     section name, and fix "Withdraw" -> "Withdrawn" ---------------------
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Carryover nominations"
$ws.Range("A9").Value  = "     Civilian, Confirmed"
$ws.Range("A10").Value = "     Civilian, Unconfirmed"
$ws.Range("A11").Value = "     Civilian, Withdrawn"

# --- Civilian lists section (rows 13-16) ----------------------------------
$ws.Range("A13").Value = "     Civilian lists, New nominations"
$ws.Range("A14").Value = "     Civilian lists, Carryover nominations"
$ws.Range("A15").Value = "     Civilian lists, Confirmed"
$ws.Range("A16").Value = "     Civilian lists, Unconfirmed"

# --- Air Force section (rows 18-20) ---------------------------------------
$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Carryover nominations"
$ws.Range("A20").Value = "     Air Force, Confirmed"

# --- Army section (rows 22-24) --------------------------------------------
$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("A23").Value = "     Army, Carryover nominations"
$ws.Range("A24").Value = "     Army, Confirmed"

# --- Navy section (rows 26-27) --------------------------------------------
$ws.Range("A26").Value = "     Navy, New nominations"
$ws.Range("A27").Value = "     Navy, Confirmed"

# --- Marine Corps section (rows 29-31) ------------------------------------
$ws.Range("A29").Value = "     Marine Corps, New nominations"
$ws.Range("A30").Value = "     Marine Corps, Carryover nominations"
$ws.Range("A31").Value = "     Marine Corps, Confirmed"

# --- Summary section: rename "Total new nominations received" and then
#     drop the now-redundant blank "Summary" header row (row 32), which
#     shifts the summary totals up by one row. ------------------------
$ws.Range("A33").Value = "Total new nominations"
$ws.Rows(32).Delete()

[void]$ws.Range("A1").Select()
